$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3500.2856
$ws.Range("I40").Value = 6001
$ws.Range("K40").Value = 6001
$ws.Range("M40").Value = -5826
$ws.Range("H62").Value = 41578.93
$ws.Range("I62").Value = 2733
$ws.Range("J62").Value = 48053.25
$ws.Range("K62").Value = 2733
$ws.Range("L62").Value = 48053.25
$ws.Range("M62").Value = -2109
$ws.Range("N62").Value = -49301.25
$ws.Range("H65").Value = 41578.93
$ws.Range("I65").Value = 2733
$ws.Range("J65").Value = 48053.25
$ws.Range("K65").Value = 13665
$ws.Range("L65").Value = 240266.25
$ws.Range("M65").Value = -10545
$ws.Range("N65").Value = -246506.25
$ws.Range("H86").Value = 3927497
$ws.Range("I86").Value = 3928.9443
$ws.Range("K86").Value = 3928.9443
$ws.Range("M86").Value = -2805.9443
$ws.Range("H89").Value = 3927497
$ws.Range("I89").Value = 3928.9443
$ws.Range("K89").Value = 19644.7215
$ws.Range("M89").Value = -14028.7215
$ws.Range("H106").Value = 6131.2856
$ws.Range("I106").Value = 1736.5
$ws.Range("K106").Value = 1736.5
$ws.Range("M106").Value = -1105.5
$ws.Range("H107").Value = 320.5
$ws.Range("I107").Value = 326.9091
$ws.Range("K107").Value = 326.9091
$ws.Range("M107").Value = 1593.0909
$ws.Range("H131").Value = 16076.846
$ws.Range("J131").Value = 14916.583
$ws.Range("L131").Value = 44749.749
$ws.Range("N131").Value = -54829.749
$ws.Range("H132").Value = 3373.3809
$ws.Range("I132").Value = 2929.575
$ws.Range("K132").Value = 8788.724999999999
$ws.Range("M132").Value = -6258.724999999999
$ws.Range("H141").Value = 11003.4

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3861.5483
$ws.Range("I32").Value = 3861.5483
$ws.Range("K32").Value = 3861.5483
$ws.Range("M32").Value = -3574.5483
$ws.Range("H61").Value = 12823869
$ws.Range("I61").Value = 23812292
$ws.Range("K61").Value = 23812292
$ws.Range("M61").Value = -23812080
$ws.Range("H88").Value = 10418382
$ws.Range("I88").Value = 27778594
$ws.Range("K88").Value = 27778594
$ws.Range("M88").Value = -27778188
$ws.Range("H91").Value = 10418382
$ws.Range("I91").Value = 27778594
$ws.Range("K91").Value = 27778594
$ws.Range("M91").Value = -27777190
$ws.Range("H136").Value = 12823869
$ws.Range("I136").Value = 23812292
$ws.Range("K136").Value = 71436876
$ws.Range("M136").Value = -71434326

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 3562.625
$ws.Range("I86").Value = 2750.1667
$ws.Range("K86").Value = 2750.1667
$ws.Range("M86").Value = -1627.1667
$ws.Range("H89").Value = 3562.625
$ws.Range("I89").Value = 2750.1667
$ws.Range("K89").Value = 13750.8335
$ws.Range("M89").Value = -8134.833500000001
$ws.Range("H99").Value = 858.7143
$ws.Range("I99").Value = 766.8333
$ws.Range("K99").Value = 766.8333
$ws.Range("M99").Value = 731.1667

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 4864.45
$ws.Range("I86").Value = 5124.2144
$ws.Range("K86").Value = 5124.2144
$ws.Range("M86").Value = -4001.2144
$ws.Range("H89").Value = 4864.45
$ws.Range("I89").Value = 5124.2144
$ws.Range("K89").Value = 25621.072
$ws.Range("M89").Value = -20005.072

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H23").Value = 100
$ws.Range("I23").Value = 100
$ws.Range("J23").Value = 100
$ws.Range("K23").Value = 300
$ws.Range("L23").Value = 300
$ws.Range("M23").Value = -65
$ws.Range("N23").Value = -770
$ws.Range("H129").Value = 2293.7368
$ws.Range("I129").Value = 1064.1428
$ws.Range("J129").Value = 3011
$ws.Range("K129").Value = 3192.4284
$ws.Range("L129").Value = 9033
$ws.Range("M129").Value = 1807.5716
$ws.Range("N129").Value = -19033
$ws.Range("H137").Value = 1316.3334
$ws.Range("I137").Value = 850
$ws.Range("J137").Value = 1549.5
$ws.Range("K137").Value = 2550
$ws.Range("L137").Value = 4648.5
$ws.Range("M137").Value = 2550
$ws.Range("N137").Value = -14848.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 3699.4
$ws.Range("I43").Value = 3699.4
$ws.Range("K43").Value = 3699.4
$ws.Range("M43").Value = -3548.4
$ws.Range("H46").Value = 38222.11
$ws.Range("J46").Value = 39874.875
$ws.Range("L46").Value = 39874.875
$ws.Range("N46").Value = -40186.875
$ws.Range("H57").Value = 49982.5
$ws.Range("J57").Value = 49982.5
$ws.Range("L57").Value = 49982.5
$ws.Range("N57").Value = -51622.5
$ws.Range("H70").Value = 12641.777
$ws.Range("I70").Value = 12051.923
$ws.Range("J70").Value = 14175.4
$ws.Range("K70").Value = 12051.923
$ws.Range("L70").Value = 14175.4
$ws.Range("M70").Value = -11781.923
$ws.Range("N70").Value = -14715.4
$ws.Range("H73").Value = 12641.777
$ws.Range("I73").Value = 12051.923
$ws.Range("J73").Value = 14175.4
$ws.Range("K73").Value = 12051.923
$ws.Range("L73").Value = 14175.4
$ws.Range("M73").Value = -11115.923
$ws.Range("N73").Value = -16047.4
$ws.Range("H102").Value = 2096.7188
$ws.Range("I102").Value = 1135.3334
$ws.Range("J102").Value = 3332.7856
$ws.Range("K102").Value = 1135.3334
$ws.Range("L102").Value = 3332.7856
$ws.Range("M102").Value = 486.6666
$ws.Range("N102").Value = -6576.7856
$ws.Range("H122").Value = 27780576
$ws.Range("I122").Value = 2337.3845
$ws.Range("K122").Value = 7012.1535
$ws.Range("M122").Value = -4562.1535
$ws.Range("H126").Value = 8130.35
$ws.Range("I126").Value = 10557.75
$ws.Range("K126").Value = 31673.25
$ws.Range("M126").Value = -29203.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 12823378
$ws.Range("I7").Value = 17859554
$ws.Range("J7").Value = 4020.818
$ws.Range("K7").Value = 17859554
$ws.Range("L7").Value = 4020.818
$ws.Range("M7").Value = -17859442
$ws.Range("N7").Value = -4244.818
$ws.Range("H40").Value = 3250.6428
$ws.Range("I40").Value = 2754.3333
$ws.Range("K40").Value = 2754.3333
$ws.Range("M40").Value = -2618.3333
$ws.Range("H43").Value = 5355799.5
$ws.Range("I43").Value = 3595571.5
$ws.Range("J43").Value = 6895999.5
$ws.Range("K43").Value = 3595571.5
$ws.Range("L43").Value = 6895999.5
$ws.Range("M43").Value = -3595378.5
$ws.Range("N43").Value = -6896385.5
$ws.Range("H53").Value = 30001
$ws.Range("J53").Value = 30001
$ws.Range("L53").Value = 30001
$ws.Range("N53").Value = -31037
$ws.Range("H126").Value = 12823378
$ws.Range("I126").Value = 17859554
$ws.Range("J126").Value = 4020.818
$ws.Range("K126").Value = 53578662
$ws.Range("L126").Value = 12062.454
$ws.Range("M126").Value = -53576192
$ws.Range("N126").Value = -17002.454
$ws.Range("H136").Value = 2156.484
$ws.Range("I136").Value = 1905.4286
$ws.Range("J136").Value = 4499.6665
$ws.Range("K136").Value = 5716.2858
$ws.Range("L136").Value = 13498.9995
$ws.Range("M136").Value = -3166.2858
$ws.Range("N136").Value = -18598.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 25051
$ws.Range("J64").Value = 0
$ws.Range("L64").Value = 0
$ws.Range("N64").ClearContents()
$ws.Range("H67").Value = 25051
$ws.Range("J67").Value = 0
$ws.Range("L67").Value = 0
$ws.Range("N67").ClearContents()
$ws.Range("H136").Value = 3956.889
$ws.Range("I136").Value = 1628.2
$ws.Range("J136").Value = 6867.75
$ws.Range("K136").Value = 4884.6
$ws.Range("L136").Value = 20603.25
$ws.Range("M136").Value = -2334.6
$ws.Range("N136").Value = -25703.25
